$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows just above the current "Total" row (row 13), pushing it down to row 15.
$ws.Rows.Item(13).Resize(2).Insert()

# Row 12 is a blank spacer row but keeps the date-column formatting (like B11).
$ws.Range("B11").Copy()
$ws.Range("B12:B13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Row 13 Hours cell picks up the same formatting as the Hours column above it.
$ws.Range("C11").Copy()
$ws.Range("C13").PasteSpecial(-4122)       # xlPasteFormats
$ws.Application.CutCopyMode = $false

# New timesheet entry in row 13.
$ws.Range("A13").Value = "Added Model classes Graphics, Collections"
$ws.Range("B13").Value = (Get-Date -Year 2018 -Month 9 -Day 14).Date
$ws.Range("C13").Value = 1.5

# The Total row (now row 15) needs its SUM formula extended to cover the new rows.
$ws.Range("C15").Formula = "=SUM(C2:C14)"

# Update the view: scroll so row 4 is at the top, and select C14.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C14").Select()
